# #5: fund, bonds, otherbonds, antique done
# Target sheet: 具有相當價值之財產 (the "valuable property" worksheet)
#
# - Row 1 (B1:E1) currently holds literal sample values instead of the
#   canonical column headers used by every other sheet; fix them up and
#   extend the header row with the standard metadata columns
#   (property_category, category, date, legislator_name, legislator_id,
#   source_file, index) in F1:L1.
# - Rows 2-3 get the matching metadata values filled into F:L.
# - Row 3's name had a stray "■" glyph in the source text; clean it up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("具有相當價值之財產")

# --- Fix the header row (B1:E1) -> canonical field names ------------------
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "quantity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"

# --- Extend header row with the standard metadata columns (F1:L1) ---------
# Copy the header style (bold + border, same as B1:E1) onto the new cells
# before writing their values.
$ws.Range("B1").Copy()
$ws.Range("F1:L1").PasteSpecial(-4122)

$ws.Range("F1").Value = "property_category"
$ws.Range("G1").Value = "category"
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"
$ws.Range("K1").Value = "source_file"
$ws.Range("L1").Value = "index"

# --- Data rows: copy the data-row style (B2) onto the new F:L cells -------
$ws.Range("B2").Copy()
$ws.Range("F2:L3").PasteSpecial(-4122)

# Row 2 (record #110 - 手錶)
$ws.Range("F2").Value = "otherbonds"
$ws.Range("G2").Value = "normal"
$ws.Range("I2").Value = "丁守中"
$ws.Range("J2").Value = 515
$ws.Range("K2").Value = "tmp8fef1"
$ws.Range("L2").Value = 110

# Row 3 (record #111 - 手錶珠寶)
$ws.Range("F3").Value = "otherbonds"
$ws.Range("G3").Value = "normal"
$ws.Range("I3").Value = "丁守中"
$ws.Range("J3").Value = 515
$ws.Range("K3").Value = "tmp8fef1"
$ws.Range("L3").Value = 111

# The "date" column (2011-11-22) looks like a date, so Excel would silently
# convert a plain .Value assignment into a date serial number. Stage it as
# text in a scratch cell first, then copy just the (already-text) value
# across so H2/H3 keep their string type instead of becoming a date.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "2011-11-22"
$ws.Range("Z1").Copy()
$ws.Range("H2").PasteSpecial(-4163)
$ws.Range("H3").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

# --- Clean up the stray "■" glyph in B3 ------------------------------------
$ws.Range("B3").Value = "手錶珠寶"
